$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date
$ws.Range("G4").Value = 44042.80167855266

# Customer name
$ws.Range("G7").Value = "Camisas Ramiro"

# Address comment box
$ws.Range("F10").Value = "qw"

# 1. Offline Storage - Digital (GB) Qty
$ws.Range("F19").Value = 240

# 2. Online Storage (GB) - piqlConnect
$ws.Range("G21").Value = 8880
$ws.Range("H21").Value = 8880

# Online Storage (GB) - Payment type / qty / unit price / total
$ws.Range("E22").Value = "yearly"
$ws.Range("F22").Value = 2000
$ws.Range("G22").Value = 0.576
$ws.Range("H22").Value = 576

# 6. Shipment cost - Reels / total
$ws.Range("E32").Value = 2
$ws.Range("H32").Value = 60

# TOTAL
$ws.Range("H33").Value = 11256

# Total to pay from the second term
$ws.Range("H34").Value = 9456
